# Tuntikirjanpito.xlsx update - "almost finished BEM + Sass transition"
#
# 1. Insert 3 blank rows before the summary block (old rows 112-114) so that
#    the summary rows move down to 115-117 while leaving rows 110-111 blank
#    to receive the two new work-log entries.
# 2. Fill rows 110 and 111 with the new data (date, hours, description, client)
#    copying the date cell's number format/alignment from an existing date row.
# 3. Re-point the hours-sum formula to include the two new rows.
# 4. Update the view's selection / scroll position to match where the user
#    was working after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the summary block (tunnit yht. / target / suoritettu(%)), currently on
# rows 112-114, down by three rows so it lands on 115-117, leaving rows
# 110-111 free for the new entries.
$ws.Range("A112:A114").EntireRow.Insert()

# New entry for 2022-01-29 ("client" work), 1 hour.
$ws.Cells.Item(104, 1).Copy()
$ws.Cells.Item(110, 1).PasteSpecial(-4122)  # xlPasteFormats - reuse the date style
$ws.Cells.Item(110, 1).Value = 44590
$ws.Cells.Item(110, 2).Value = 1
$ws.Cells.Item(110, 3).Value = "Notification scss, pieniä muutoksia, LoginBar scss"
$ws.Cells.Item(110, 4).Value = "client"

# Second new entry, same day, 3 hours.
$ws.Cells.Item(111, 2).Value = 3
$ws.Cells.Item(111, 3).Value = "Reset, Base, Header, LoginBar loppuun, Nav, Footer, Loading, Authenticate   scss, mixinien luontia ja kokeilua"
$ws.Cells.Item(111, 4).Value = "client"

# Extend the hours-total formula to cover the two new rows.
$ws.Cells.Item(115, 2).Formula = "=SUM(B2:B111)"

# Match the saved view state (scrolled down a bit further, selection moved).
$ws.Range("C111").Select()
$excel.ActiveWindow.ScrollRow = 100
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Left = 30975
$excel.ActiveWindow.Top = 2205
